# CAD-1156 add the recommitment columns
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New header labels to add after column AL (38) -> AM(39..42)
$newHeaders = @("recommitment", "recommitment start date", "recommitment end date", "external reference id")

$startCol = 39  # AM
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = $startCol + $i
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $newHeaders[$i]
}

# Set column widths to match target (bestFit/customWidth values from diff)
$ws.Columns.Item(39).ColumnWidth = 12.6640625
$ws.Columns.Item(40).ColumnWidth = 20.33203125
$ws.Columns.Item(41).ColumnWidth = 19.6640625
$ws.Columns.Item(42).ColumnWidth = 16.6640625

# Update the view: top-left cell and selection
$ws.Application.ActiveWindow.ScrollColumn = 25  # Y column
$ws.Range("AL15").Select()
